$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new bold paragraph "mod_quantitation.R / fct_quantitation.R"
#    right before the "mod_derived_traits.R / fct_derived_traits.R" paragraph
#    (exclude quantitation clusters from spectra curation).
# ---------------------------------------------------------------------------

$targetRange = $word.ActiveDocument.Content
$targetRange.Find.Execute("mod_derived_traits.R", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hit = $targetRange.Duplicate
$hit.Expand(4) | Out-Null
$derivedPara = $hit.Paragraphs(1)

$derivedPara.Range.InsertParagraphBefore()
$newRange = $derivedPara.Range

$quantXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val="4472C4" w:themeColor="accent1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>mod_quantitation.R</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> / </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:color w:val="4472C4" w:themeColor="accent1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>fct_quantitation.R</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@

$newRange.InsertXML($quantXml)

# ---------------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from the "mod_data_exploration.R" run to
#    the tab run that immediately precedes it (in the
#    "mod_tab_repeatability_plot.R" paragraph).
# ---------------------------------------------------------------------------

$d2 = $word.ActiveDocument

$findRange = $d2.Content
$findRange.Find.Execute("mod_tab_repeatability_plot.R", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hit2 = $findRange.Duplicate
$hit2.Expand(4) | Out-Null
$tabPlotPara = $hit2.Paragraphs(1)

$tabPlotXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>mod_tab_repeatability_plot.R</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@
$tabPlotPara.Range.InsertXML($tabPlotXml)

$findRange2 = $d2.Content
$findRange2.Find.Execute("mod_data_exploration.R", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hit3 = $findRange2.Duplicate
$hit3.Expand(4) | Out-Null
$dataExpPara = $hit3.Paragraphs(1)

$dataExpXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>mod_data_exploration.R</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> / </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:color w:val="4472C4" w:themeColor="accent1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>fct_data_exploration.R</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
'@
$dataExpPara.Range.InsertXML($dataExpXml)

Write-Output "edit complete"
